$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: "Renderer-independent blend mode" now has a "Yes" marker in column C,
# styled the same as the other "Yes" cells (e.g. B7) with the green fill.
$ws.Range("B7").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Value = "Yes"

# Row 26: "Input Method API" now has a "Yes" marker in column B,
# styled the same way.
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null
$ws.Range("B26").Value = "Yes"

# Clear clipboard marching ants / pasted-range state.
$excel.CutCopyMode = 0

# Update the saved selection to A12 (anchor) across A12:A14.
$null = $ws.Range("A12:A14").Select()
